# Update "想去人数" (want-to-go count) figures for three events that each
# appear on multiple sheets (展览 / 演出 / 全部类型).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 7915
$ws1.Range("F15").Value = 47
$ws1.Range("F33").Value = 1900
$ws1.Range("F38").Value = 3747
$ws1.Range("F39").Value = 336
$ws1.Range("F40").Value = 284
$ws1.Range("F44").Value = 2

# --- Sheet "演出" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 142

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 7915
$ws4.Range("F18").Value = 47
$ws4.Range("F36").Value = 1900
$ws4.Range("F41").Value = 336
$ws4.Range("F42").Value = 284
$ws4.Range("F47").Value = 142

$wb.Save()
